$d = $word.ActiveDocument

# --- 1. Honor-code quote: merge the two runs that were split around the
#    old cursor / _GoBack bookmark back into a single run. ---
$d.Content.Find.Execute(
    "nor will I accept the actions of those who do.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "nor will I accept the actions of those who do.",
    2
) | Out-Null

# The old edit-point bookmark lived between those two runs; it no longer
# belongs there once the runs are merged (Word will replant it at the
# newest edit location below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Rename the "Grievances" heading to "Course Concerns". ---
$d.Content.Find.Execute(
    "Grievances",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Course Concerns",
    2
) | Out-Null

# --- 3. Replant _GoBack right after "Course Concerns" (end of that
#    paragraph's run, before the paragraph mark) -- this is where Word
#    leaves the mark after the most recent edit. ---
$findRng = $d.Content
$findRng.Find.Execute(
    "Course Concerns",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "",
    0
) | Out-Null
$endPos = $findRng.End

# Adding a bookmark collapsed exactly at "end of last run in paragraph"
# needs a non-boundary anchor, so insert a throwaway marker character
# right after the text, anchor the bookmark to the position just before
# it, then delete the marker -- the bookmark stays put.
$markerRng = $d.Range($endPos, $endPos)
$markerRng.InsertAfter("X")

$bmRng = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRng) | Out-Null

$d.Range($endPos, $endPos + 1).Delete()
